# Form the consolidated report: fix the "Absent" (column H) values so
# that they correctly reflect whether the student attended that day
# (column D, "Total Attendance Count").
#
# Rows 6, 12, 15, 19 were incorrectly marked as 0 (present) even though
# the student did not attend that day -> they must be 1 (absent).
#
# Rows 7, 13, 16, 20 were left blank (inline string placeholder) even
# though the student did attend that day -> they must be filled in with
# the numeric value 0 (not absent).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0

$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0

$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 0

$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
